# daily auto push: 2026-02-02 10:03 UTC
#
# The sheet contains a long time-series table (date / weekday / hour / rank)
# that runs from row 2 through row 800 (row 1 is the header).
# This "daily auto push" inserts one new data point
#   (2026/02/02, 月, 16, 201)
# right before the existing 2026/12/29 block (row 759), pushing every
# subsequent row down by one. The sheet's used range therefore grows from
# A1:D800 to A1:D801.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 759; Excel automatically shifts rows
# 759:800 down to 760:801 and extends the sheet dimension accordingly.
$ws.Rows.Item(759).Insert()

# Fill in the newly inserted row with the new measurement.
$dateCell = $ws.Cells.Item(759, 1)

# Force the date column to be written as plain text (matching every other
# row in the column, which stores dates as literal "yyyy/mm/dd" strings
# rather than real date serials) instead of letting Excel auto-convert the
# "2026/02/02" string into a date value.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/02"
# Drop the temporary text formatting again so the cell ends up with the
# same (default) style as its neighbours.
$dateCell.ClearFormats()

$ws.Cells.Item(759, 2).Value = "月"
$ws.Cells.Item(759, 3).Value = 16
$ws.Cells.Item(759, 4).Value = 201
